$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# Row 11
$ws.Range("H11").Value = 62.6
$ws.Range("I11").Value = 62.6
$ws.Range("K11").Value = 62.6
$ws.Range("M11").Value = 77.40000000000001

# Row 98
$ws.Range("H98").Value = 2233.6667
$ws.Range("I98").Value = 2407.6155
$ws.Range("J98").Value = 1103
$ws.Range("K98").Value = 2407.6155
$ws.Range("L98").Value = 1103
$ws.Range("M98").Value = -909.6154999999999
$ws.Range("N98").Value = -4099

# Row 113
$ws.Range("H113").Value = 5576.375
$ws.Range("I113").Value = 5121
$ws.Range("J113").Value = 6335.3335
$ws.Range("K113").Value = 5121
$ws.Range("L113").Value = 6335.3335
$ws.Range("M113").Value = -1867
$ws.Range("N113").Value = -12843.3335

# Row 122
$ws.Range("H122").Value = 2233.6667
$ws.Range("I122").Value = 2407.6155
$ws.Range("J122").Value = 1103
$ws.Range("K122").Value = 7222.8465
$ws.Range("L122").Value = 3309
$ws.Range("M122").Value = -4772.8465
$ws.Range("N122").Value = -8209

# Row 132
$ws.Range("H132").Value = 2440973.2
$ws.Range("I132").Value = 1970
$ws.Range("J132").Value = 12501861
$ws.Range("K132").Value = 5910
$ws.Range("L132").Value = 37505583
$ws.Range("M132").Value = -3380
$ws.Range("N132").Value = -37510643

# Row 141
$ws.Range("H141").Value = 1836.6666
$ws.Range("I141").Value = 1836.6666
$ws.Range("K141").Value = 5509.9998
$ws.Range("M141").Value = -329.9997999999996

$ws = $wb.Worksheets("ARM")
# Row 45
$ws.Range("H45").Value = 974.0769
$ws.Range("I45").Value = 843.1111
$ws.Range("J45").Value = 1268.75
$ws.Range("K45").Value = 843.1111
$ws.Range("L45").Value = 1268.75
$ws.Range("M45").Value = -466.1111
$ws.Range("N45").Value = -2022.75

# Row 61
$ws.Range("H61").Value = 1432.1464
$ws.Range("I61").Value = 1473.081
$ws.Range("J61").Value = 1053.5
$ws.Range("K61").Value = 1473.081
$ws.Range("L61").Value = 1053.5
$ws.Range("M61").Value = -1261.081
$ws.Range("N61").Value = -1477.5

# Row 63
$ws.Range("H63").Value = 2217.077
$ws.Range("I63").Value = 2135.6667
$ws.Range("K63").Value = 2135.6667
$ws.Range("M63").Value = -1449.6667

# Row 66
$ws.Range("H66").Value = 2217.077
$ws.Range("I66").Value = 2135.6667
$ws.Range("K66").Value = 10678.3335
$ws.Range("M66").Value = -7246.333500000001

# Row 74
$ws.Range("H74").Value = 1373.5
$ws.Range("I74").Value = 1389.56
$ws.Range("J74").Value = 1293.2
$ws.Range("K74").Value = 1389.56
$ws.Range("L74").Value = 1293.2
$ws.Range("M74").Value = -515.5599999999999
$ws.Range("N74").Value = -3041.2

# Row 77
$ws.Range("H77").Value = 1373.5
$ws.Range("I77").Value = 1389.56
$ws.Range("J77").Value = 1293.2
$ws.Range("K77").Value = 6947.799999999999
$ws.Range("L77").Value = 6466
$ws.Range("M77").Value = -2579.799999999999
$ws.Range("N77").Value = -15202

# Row 102
$ws.Range("H102").Value = 3010
$ws.Range("I102").Value = 3010
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3010
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1388
$ws.Range("N102").ClearContents()

# Row 122
$ws.Range("H122").Value = 907.5833
$ws.Range("I122").Value = 871.9091
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 2615.7273
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -165.7273
$ws.Range("N122").Value = -8800

# Row 132
$ws.Range("H132").Value = 1518.8
$ws.Range("I132").Value = 1319.3334
$ws.Range("K132").Value = 3958.0002
$ws.Range("M132").Value = -1428.0002

# Row 136
$ws.Range("H136").Value = 1432.1464
$ws.Range("I136").Value = 1473.081
$ws.Range("J136").Value = 1053.5
$ws.Range("K136").Value = 4419.242999999999
$ws.Range("L136").Value = 3160.5
$ws.Range("M136").Value = -1869.242999999999
$ws.Range("N136").Value = -8260.5

$ws = $wb.Worksheets("BSM")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

# Row 134
$ws.Range("H134").Value = 2130.946
$ws.Range("I134").Value = 1838.9062
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 5516.7186
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -2981.7186
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets("CRP")
# Row 99
$ws.Range("H99").Value = 3259.7036
$ws.Range("I99").Value = 2277.8462
$ws.Range("J99").Value = 4171.4287
$ws.Range("K99").Value = 2277.8462
$ws.Range("L99").Value = 4171.4287
$ws.Range("M99").Value = -779.8462
$ws.Range("N99").Value = -7167.4287

# Row 126
$ws.Range("H126").Value = 3259.7036
$ws.Range("I126").Value = 2277.8462
$ws.Range("J126").Value = 4171.4287
$ws.Range("K126").Value = 6833.5386
$ws.Range("L126").Value = 12514.2861
$ws.Range("M126").Value = -4363.5386
$ws.Range("N126").Value = -17454.2861

$ws = $wb.Worksheets("CUL")
# Row 102
$ws.Range("H102").Value = 3833.3333
$ws.Range("J102").Value = 3833.3333
$ws.Range("L102").Value = 11499.9999
$ws.Range("N102").Value = -16367.9999

# Row 113
$ws.Range("H113").Value = 664.73914
$ws.Range("I113").Value = 786.6667
$ws.Range("J113").Value = 621.7059
$ws.Range("K113").Value = 2360.0001
$ws.Range("L113").Value = 1865.1177
$ws.Range("M113").Value = -190.0001000000002
$ws.Range("N113").Value = -6205.117700000001

# Row 124
$ws.Range("H124").Value = 5000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 5000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 15000
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -24820

# Row 125
$ws.Range("H125").Value = 3088.3333
$ws.Range("I125").Value = 2015
$ws.Range("J125").Value = 3625
$ws.Range("K125").Value = 6045
$ws.Range("L125").Value = 10875
$ws.Range("M125").Value = -1125
$ws.Range("N125").Value = -20715

# Row 131
$ws.Range("H131").Value = 14765932
$ws.Range("I131").Value = 62625276
$ws.Range("J131").Value = 39979.77
$ws.Range("K131").Value = 187875828
$ws.Range("L131").Value = 119939.31
$ws.Range("M131").Value = -187870788
$ws.Range("N131").Value = -130019.31

$ws = $wb.Worksheets("GSM")
# Row 70
$ws.Range("H70").Value = 4397.923
$ws.Range("I70").Value = 4260.933
$ws.Range("K70").Value = 4260.933
$ws.Range("M70").Value = -3990.933

# Row 73
$ws.Range("H73").Value = 4397.923
$ws.Range("I73").Value = 4260.933
$ws.Range("K73").Value = 4260.933
$ws.Range("M73").Value = -3324.933

# Row 102
$ws.Range("H102").Value = 1675.3334
$ws.Range("I102").Value = 1516.5
$ws.Range("J102").Value = 1993
$ws.Range("K102").Value = 1516.5
$ws.Range("L102").Value = 1993
$ws.Range("M102").Value = 105.5
$ws.Range("N102").Value = -5237

# Row 126
$ws.Range("H126").Value = 2512.3333
$ws.Range("I126").Value = 2690.0833
$ws.Range("J126").Value = 1801.3334
$ws.Range("K126").Value = 8070.249899999999
$ws.Range("L126").Value = 5404.0002
$ws.Range("M126").Value = -5600.249899999999
$ws.Range("N126").Value = -10344.0002

$ws = $wb.Worksheets("LTW")
# Row 7
$ws.Range("H7").Value = 2297.6365
$ws.Range("I7").Value = 2297.6365
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2297.6365
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2185.6365
$ws.Range("N7").ClearContents()

# Row 40
$ws.Range("H40").Value = 1444732
$ws.Range("I40").Value = 3367904.8
$ws.Range("J40").Value = 2352.5
$ws.Range("K40").Value = 3367904.8
$ws.Range("L40").Value = 2352.5
$ws.Range("M40").Value = -3367768.8
$ws.Range("N40").Value = -2624.5

# Row 46
$ws.Range("H46").Value = 2500
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2312
$ws.Range("N46").ClearContents()

# Row 126
$ws.Range("H126").Value = 2297.6365
$ws.Range("I126").Value = 2297.6365
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6892.9095
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4422.9095
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets("WVR")
# Row 81
$ws.Range("H81").Value = 2072.611
$ws.Range("I81").Value = 2025.5
$ws.Range("J81").Value = 2237.5
$ws.Range("K81").Value = 4051
$ws.Range("L81").Value = 4475
$ws.Range("M81").Value = -2990
$ws.Range("N81").Value = -6597

# Row 84
$ws.Range("H84").Value = 2072.611
$ws.Range("I84").Value = 2025.5
$ws.Range("J84").Value = 2237.5
$ws.Range("K84").Value = 20255
$ws.Range("L84").Value = 22375
$ws.Range("M84").Value = -14951
$ws.Range("N84").Value = -32983

# Row 126
$ws.Range("H126").Value = 916.8148
$ws.Range("I126").Value = 857.8
$ws.Range("J126").Value = 1085.4286
$ws.Range("K126").Value = 2573.4
$ws.Range("L126").Value = 3256.2858
$ws.Range("M126").Value = -103.3999999999996
$ws.Range("N126").Value = -8196.2858

# Row 130
$ws.Range("H130").Value = 36266.668
$ws.Range("J130").Value = 36266.668
$ws.Range("L130").Value = 36266.668
$ws.Range("N130").Value = -46306.668

# Row 132
$ws.Range("H132").Value = 753.5714
$ws.Range("I132").Value = 627.63635
$ws.Range("J132").Value = 1861.8
$ws.Range("K132").Value = 1882.90905
$ws.Range("L132").Value = 5585.4
$ws.Range("M132").Value = 647.09095
$ws.Range("N132").Value = -10645.4
